$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-9 (MuSCs -> ECs/Resolving-Mac rows removed entirely)
$ws.Range("A6:T9").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt7b"
$ws.Range("C2").Value = "Fzd10"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.021087
$ws.Range("H2").Value = 0.063261
$ws.Range("I2").Value = 0.02328126719340038
$ws.Range("J2").Value = 0.02328126719340038
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.085107
$ws.Range("N2").Value = 0.255321
$ws.Range("O2").Value = 0.803017436594203
$ws.Range("P2").Value = 0.8030174365942029
$ws.Range("Q2").Value = 0.001794651309
$ws.Range("R2").Value = 0.016151861781
$ws.Range("S2").Value = 0.01869526350230909
$ws.Range("T2").Value = 0.01869526350230909

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt7b"
$ws.Range("C3").Value = "Fzd10"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.021087
$ws.Range("H3").Value = 0.063261
$ws.Range("I3").Value = 0.02328126719340038
$ws.Range("J3").Value = 0.02328126719340038
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.020877
$ws.Range("N3").Value = 0.06263099999999999
$ws.Range("O3").Value = 0.1969825634057971
$ws.Range("P3").Value = 0.1969825634057971
$ws.Range("Q3").Value = 0.0004402332989999999
$ws.Range("R3").Value = 0.003962099691
$ws.Range("S3").Value = 0.004586003691091294
$ws.Range("T3").Value = 0.004586003691091294

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Wnt7b"
$ws.Range("C4").Value = "Fzd10"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8846626666666667
$ws.Range("H4").Value = 2.653988
$ws.Range("I4").Value = 0.9767187328065996
$ws.Range("J4").Value = 0.9767187328065997
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.085107
$ws.Range("N4").Value = 0.255321
$ws.Range("O4").Value = 0.803017436594203
$ws.Range("P4").Value = 0.8030174365942029
$ws.Range("Q4").Value = 0.07529098557200001
$ws.Range("R4").Value = 0.677618870148
$ws.Range("S4").Value = 0.7843221730918939
$ws.Range("T4").Value = 0.7843221730918939

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Wnt7b"
$ws.Range("C5").Value = "Fzd10"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.8846626666666667
$ws.Range("H5").Value = 2.653988
$ws.Range("I5").Value = 0.9767187328065996
$ws.Range("J5").Value = 0.9767187328065997
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.020877
$ws.Range("N5").Value = 0.06263099999999999
$ws.Range("O5").Value = 0.1969825634057971
$ws.Range("P5").Value = 0.1969825634057971
$ws.Range("Q5").Value = 0.018469102492
$ws.Range("R5").Value = 0.166221922428
$ws.Range("S5").Value = 0.1923965597147058
$ws.Range("T5").Value = 0.1923965597147058
